$wb = $excel.ActiveWorkbook

# --- Sheet "Programacao": add new row 26 ---
$ws1 = $wb.Worksheets.Item("Programacao")
$ws1.Range("A26").Value = "13/02/2025"
$ws1.Range("B26").Value = "15:54"
$ws1.Range("C26").Value = "JACO PEREIRA DANTES"
$ws1.Range("D26").Value = "(84)996589874"
$ws1.Range("E26").Value = "355890"
$ws1.Range("F26").Value = "NORSAL"
$ws1.Range("G26").Value = 29700
$ws1.Range("H26").Value = "RYE3J61"
$ws1.Range("I26").Value = "BAU"
$ws1.Range("J26").Value = "SAL REFINADO 25 KG"
$ws1.Range("K26").Value = "MONTE SERENO"

# --- Sheet "Planilha": add new row 41 ---
$ws2 = $wb.Worksheets.Item("Planilha")
$ws2.Range("A41").Value = "ENTRADA"
$ws2.Range("B41").Value = "13/02/2025"
$ws2.Range("C41").Value = "RYE3J61"
$ws2.Range("D41").Value = "BAU"
$ws2.Range("E41").Value = "MONTE SERENO"
$ws2.Range("F41").Value = "SAL REFINADO"
$ws2.Range("G41").Value = "25 KG"
$ws2.Range("H41").Value = "NORSAL"
$ws2.Range("I41").Value = "355890"
$ws2.Range("J41").Value = "355889"
$ws2.Range("K41").Value = 22
$ws2.Range("L41").Value = "140303425"
$ws2.Range("M41").Value = "fev/27"
$ws2.Range("N41").Value = 29700

# --- Sheet "Descarga do Sal": update form with latest entry data ---
$ws3 = $wb.Worksheets.Item("Descarga do Sal")
$ws3.Range("D8").Value = "13/02/2025"
$ws3.Range("K8").Value = "15:54"
$ws3.Range("D10").Value = "JACO PEREIRA DANTES"
$ws3.Range("D12").Value = "(84)996589874"
$ws3.Range("D14").Value = "BAU"
$ws3.Range("K14").Value = "RYE3J61"
$ws3.Range("D16").Value = "MONTE SERENO"
$ws3.Range("D20").Value = "355890"
$ws3.Range("K20").Value = "355889"
$ws3.Range("P20").Value = 29700
$ws3.Range("D22").Value = ""
$ws3.Range("K22").Value = ""
$ws3.Range("P22").Value = ""
$ws3.Range("D24").Value = ""
$ws3.Range("K24").Value = ""
$ws3.Range("P24").Value = ""
$ws3.Range("D26").Value = "SAL REFINADO"
$ws3.Range("L26").Value = "fev/27"
$ws3.Range("D28").Value = "140303425"
$ws3.Range("H28").Value = "355890"
$ws3.Range("K28").Value = 29700
$ws3.Range("O28").Value = 22
$ws3.Range("D30").Value = ""
$ws3.Range("H30").Value = ""
$ws3.Range("K30").Value = ""
$ws3.Range("O30").Value = ""
$ws3.Range("D32").Value = ""
$ws3.Range("H32").Value = ""
$ws3.Range("K32").Value = ""
$ws3.Range("O32").Value = ""
